# #5: property boat&car done
#
# The "汽車" (car) sheet's header row was accidentally populated with the
# row-2 data values instead of real column headers, and the sheet was
# missing the "capacity" column plus the trailing metadata columns
# (property_category / category / date / legislator_name / legislator_id /
# source_file / index) that every other property sheet already carries.
# Fix the headers and fill in the missing columns on the data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("汽車")

# --- Row 1: real header labels (was holding stray data) -----------------
$ws.Cells.Item(1, 2).Value  = "name"
$ws.Cells.Item(1, 3).Value  = "capacity"
$ws.Cells.Item(1, 4).Value  = "owner"
$ws.Cells.Item(1, 5).Value  = "register_date"
$ws.Cells.Item(1, 6).Value  = "register_reason"
$ws.Cells.Item(1, 7).Value  = "acquire_value"
$ws.Cells.Item(1, 8).Value  = "property_category"
$ws.Cells.Item(1, 9).Value  = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# Give the new header cells the same bold/centered/bordered look as B1:G1.
$ws.Range("G1").Copy()
$ws.Range("H1:N1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 2: keep the existing data, append the missing trailing columns -
$ws.Cells.Item(2, 8).Value  = "land"
$ws.Cells.Item(2, 9).Value  = "normal"

# "date" is a plain text field (yyyy-mm-dd string), not a real Excel date -
# force text formatting first so the COM layer doesn't coerce it into a
# date serial number, then restore the plain formatting used elsewhere.
$ws.Cells.Item(2, 10).NumberFormat = "@"
$ws.Cells.Item(2, 10).Value = "2011-11-22"
$ws.Range("G2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(2, 11).Value = "丁守中"
$ws.Cells.Item(2, 12).Value = 515
$ws.Cells.Item(2, 13).Value = "tmp8fef1"
$ws.Cells.Item(2, 14).Value = 46

# Match the plain (unbolded, borderless) data-row formatting used on B2:G2.
$ws.Range("G2").Copy()
$ws.Range("H2:I2").PasteSpecial(-4122)
$ws.Range("K2:N2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
